$wb = $excel.ActiveWorkbook

# ---- Sheet: Resource Overview ----
$ws = $wb.Worksheets.Item("Resource Overview")
$ws.Range("A2").Value = "Product Implementation Project"
$ws.Range("B6").Value = "Enterprise Product Implementation"
$ws.Range("A18").Value = "Product Design/Product"
$ws.Range("G18").Value = "Development, Python, Statistics"
$ws.Range("A20").Value = "Manufacturing Engineering"
$ws.Range("A22").Value = "Production Operations/Infrastructure"
$ws.Range("G23").Value = "Manufacturing, Communication"
# row 4 already exists blank in the source; keep it materialized through the round trip
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(13, 1).Style = "Normal"

# ---- Sheet: Detailed Staffing Plan ----
$ws = $wb.Worksheets.Item("Detailed Staffing Plan")
$ws.Range("A1").Value = "DETProductLED STAFFING PLAN"
$ws.Range("B9").Value = "Lead Product Designer"
$ws.Range("C9").Value = "Product Design/Product"
$ws.Range("K9").Value = "Development, Advanced Engineering, Python"
$ws.Range("P9").Value = "Product Lead"
$ws.Range("B10").Value = "Senior Product Designer"
$ws.Range("C10").Value = "Product Design/Product"
$ws.Range("K10").Value = "Development, Statistics, R/Python"
$ws.Range("B11").Value = "Product Designer"
$ws.Range("C11").Value = "Product Design/Product"
$ws.Range("K11").Value = "Development, Python, Visualization"
$ws.Range("B12").Value = "Development Engineer"
$ws.Range("C12").Value = "Product Design/Product"
$ws.Range("K12").Value = "DevelopmentOps, Python, Cloud"
$ws.Range("B13").Value = "Junior Product Designer"
$ws.Range("C13").Value = "Product Design/Product"
$ws.Range("B18").Value = "Senior Manufacturing Engineer"
$ws.Range("C18").Value = "Manufacturing Engineering"
$ws.Range("K18").Value = "ETL, Spark, Manufacturing Systems"
$ws.Range("B19").Value = "Manufacturing Engineer"
$ws.Range("C19").Value = "Manufacturing Engineering"
$ws.Range("K19").Value = "SQL, Python, Data Production Lines"
$ws.Range("B20").Value = "Cloud Manufacturing Engineer"
$ws.Range("C20").Value = "Manufacturing Engineering"
$ws.Range("B23").Value = "Production Operations Engineer"
$ws.Range("C23").Value = "Production Operations/Infrastructure"
$ws.Range("P23").Value = "Production Operations Lead"
$ws.Range("C24").Value = "Production Operations/Infrastructure"
$ws.Range("K25").Value = "Change Management, Manufacturing"
$ws.Range("B26").Value = "Manufacturing Specialist"
$ws.Range("K26").Value = "Manufacturing Design, Facilitation"
$ws.Cells.Item(2, 1).Style = "Normal"

# ---- Sheet: Resource Timeline ----
$ws = $wb.Worksheets.Item("Resource Timeline")
$ws.Range("B5").Value = "Lead Product Designer"
$ws.Range("B7").Value = "Senior Manufacturing Engineer"
$ws.Range("B9").Value = "Production Operations Engineer"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(11, 1).Style = "Normal"

# ---- Sheet: Skills Matrix ----
$ws = $wb.Worksheets.Item("Skills Matrix")
$ws.Range("C3").Value = "CAD/Design Tools"
$ws.Range("D3").Value = "Product Engineering"
$ws.Range("E3").Value = "Manufacturing Engineering"
$ws.Range("F3").Value = "Manufacturing Systems"
$ws.Range("J3").Value = "Production Operations"
$ws.Range("B5").Value = "Lead Product Designer"
$ws.Range("B7").Value = "Senior Manufacturing Engineer"
$ws.Range("B9").Value = "Production Operations Engineer"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(11, 1).Style = "Normal"

# ---- Sheet: Cost Analysis ----
$ws = $wb.Worksheets.Item("Cost Analysis")
$ws.Range("A6").Value = "Product Design/Product"
$ws.Range("A8").Value = "Manufacturing Engineering"
$ws.Range("A10").Value = "Production Operations/Infrastructure"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(14, 1).Style = "Normal"
$ws.Cells.Item(15, 1).Style = "Normal"

# ---- Sheet: Resource Risk Assessment ----
$ws = $wb.Worksheets.Item("Resource Risk Assessment")
$ws.Range("B5").Value = "Team lacks required Development expertise"
$ws.Range("F5").Value = "Manufacturing and external consultants"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(13, 1).Style = "Normal"
